$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New TPM-based data for rows 2-4 (A:T). Row order / cluster labels change,
# rows 5-7 are removed entirely.
$data = @(
    @("FAPs","Tac1","Tacr2","ECs",3,1,1.545371333333333,4.636114,1,1,1,0.3333333333333333,0.02438133333333333,0.073144,0.2065981245057056,0.2065981245057056,0.03767821360177778,0.339103922416,0.2065981245057056,0.2065981245057056),
    @("FAPs","Tac1","Tacr2","FAPs",3,1,1.545371333333333,4.636114,1,1,1,0.3333333333333333,0.070577,0.211731,0.5980425940571686,0.5980425940571686,0.1090676725926667,0.9816090533340001,0.5980425940571686,0.5980425940571686),
    @("FAPs","Tac1","Tacr2","MuSCs",3,1,1.545371333333333,4.636114,1,1,1,0.3333333333333333,0.023055,0.069165,0.1953592814371257,0.1953592814371257,0.03562853609,0.32065682481,0.1953592814371257,0.1953592814371257)
)

# Remove old rows 5-7 (there are only 4 rows of data after this change).
$ws.Rows.Item(5).Resize(3).Delete() | Out-Null

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
